# Save last stable rank
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "last stable rank" values (Classement / Nombre de courses)
$ws.Range("C2").Value = 1037
$ws.Range("D2").Value = 13

$ws.Range("C3").Value = 998
$ws.Range("D3").Value = 11

$ws.Range("C4").Value = 1176
$ws.Range("D4").Value = 25

$ws.Range("C5").Value = 896
$ws.Range("D5").Value = 15

# Restore the view: scroll back to show column A (drop the old
# topLeftCell="O1") and move the selection to K14.
$ws.Activate()
$ws.Range("A1").Select() | Out-Null
$ws.Range("K14").Select() | Out-Null
